# Error Calculations and Plots
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are no longer present in the cleaned data set:
# "RM 232" (originally row 26) and "SC 92" (originally row 28, which becomes
# row 27 once the first deletion has shifted everything up by one).
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Fill in / clear individual cells that changed between the two
# imputed / missing-data snapshots.
$ws.Range("D3").Value = -14.2

$ws.Range("D5").Value = "'"
$ws.Range("D5").Style = "Normal"

$ws.Range("D21").Value = -14.3

$ws.Range("D23").Value = "'"
$ws.Range("D23").Style = "Normal"

$ws.Range("D32").Value = -14.7
